# groom_bout.xlsx edit script
#
# Summary of the change (see commit message):
#  - The "form_id" setting is renamed to "table_id" (same value: groom_bout)
#    on the "settings" sheet.
#  - A new "properties" worksheet is appended that will drive properties.csv:
#    columns partition/aspect/key/type/value, with a single data row
#    describing the Table-level "default" "colOrder" array property.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. survey sheet: selection lands back on C10 (unchanged location, but it
#    will stop being the tab-displayed-as-active sheet once "properties" is
#    selected last below).
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Rows.Item(10).RowHeight = 26.2
$survey.Range("C10").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. settings sheet: "form_id" -> "table_id" (the stored value, groom_bout,
#    stays the same) and the active selection moves to A3.
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "table_id"
$settings.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. New "properties" sheet, appended after "settings" and left as the
#    active sheet/tab.
# ---------------------------------------------------------------------------
$props = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$props.Name = "properties"

$props.Range("A1").Value = "partition"
$props.Range("B1").Value = "aspect"
$props.Range("C1").Value = "key"
$props.Range("D1").Value = "type"
$props.Range("E1").Value = "value"

$props.Range("A2").Value = "Table"
$props.Range("B2").Value = "default"
$props.Range("C2").Value = "colOrder"
$props.Range("D2").Value = "array"
$props.Range("E2").Value = '["GRM_FOL_date","GRM_FOL_B_focal_AnimID","GRM_time_begin","GRM_B_partner_AnimID","GRM_time_end","GRM_dur","GRM_direction","GRM_time_certainty","GRM_other_partners"]'

$props.Columns("A:E").ColumnWidth = 13

$props.Range("E11").Select() | Out-Null
